{"js": "const results = context.document.body.search(\"WORKPLACE\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nconst rng = results.items[0];\nconst rngEnd = rng.getRange('End');\nconst cmt = rngEnd.insertComment(\"Test comment collapsed end\");\nawait context.sync();\nreturn \"done\";\n", "ps1": "$d = $word.ActiveDocument\n$p1 = $d.Paragraphs(1)\n$r = $p1.Range\n[void]$r.MoveEnd(1, -1)\nWrite-Output (\"Text=\" + $r.Text)\nWrite-Output (\"Start=\" + $r.Start)\nWrite-Output (\"End=\" + $r.End)\n$r.Collapse(0)\nWrite-Output (\"CollapsedStart=\" + $r.Start)\nWrite-Output (\"CollapsedEnd=\" + $r.End)\n$d.Comments.Add($r, \"Test comment com4\")\n"}
